# initial run WC playoffs
# Append the newly opened week-19 (Wild Card) lines to the betting log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newGames = @(
    @(19, "LA_CAR",  46.5, -10),
    @(19, "GB_CHI",  45.5, -1.5),
    @(19, "BUF_JAX", 51.5, -1.5),
    @(19, "SF_PHI",  46.5,  3.5),
    @(19, "LAC_NE",  45.5,  3.5),
    @(19, "HOU_PIT", 38.5, -3.5)
)

$xlUp = -4162
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End($xlUp).Row
$startRow = $lastRow + 1

for ($i = 0; $i -lt $newGames.Length; $i++) {
    $row = $startRow + $i
    $game = $newGames[$i]
    $ws.Cells.Item($row, 1).Value = $game[0]
    $ws.Cells.Item($row, 2).Value = $game[1]
    $ws.Cells.Item($row, 3).Value = $game[2]
    $ws.Cells.Item($row, 4).Value = $game[3]
}

$lastDataRow = $startRow + $newGames.Length - 1
$ws.Cells.Item($lastDataRow + 1, 4).Select()
